# Daily attendance processing - 2025-10-15 22:44:44
# Reorders the "Recorded By" (column G) list for specific rows so that the
# two comma-separated entries swap positions (e.g. "a, b" -> "b, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(3,4,6,10,11,12,13,14,15,17,18,19,30,31,33,37,38,39,40,41,42,44,45,46,57,58,60,64,65,66,67,68,69,71,72,73,86,87,88,89,90,93,95,96,97,112,113,114,115,116,119,121,122,123,138,139,140,141,142,145,147,148,149)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = $cell.Value2
    $parts = $current -split ', ', 2
    if ($parts.Count -eq 2) {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
}
